$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update job_title (column E) for the users table rows.
# Only rows whose job title text actually changes are touched; rows 11
# (Secrétaire médicale) and 19 (Assistante de direction) keep their text.
$ws.Range("E2").Value  = "Médecin associé"
$ws.Range("E3").Value  = "Médecin associée"
$ws.Range("E4").Value  = "Médecin associé"
$ws.Range("E5").Value  = "Médecin associé"
$ws.Range("E6").Value  = "Médecin associée"
$ws.Range("E7").Value  = "Médecin associée"
$ws.Range("E8").Value  = "Médecin associé"
$ws.Range("E9").Value  = "Cofondateur & cogérant"
$ws.Range("E10").Value = "Médecin associé"
$ws.Range("E12").Value = "Médecin associé"
$ws.Range("E13").Value = "Médecin associé"
$ws.Range("E14").Value = "Médecin associé"
$ws.Range("E15").Value = "Cogérant"
$ws.Range("E16").Value = "Médecin associé"
$ws.Range("E17").Value = "Médecin associée"
$ws.Range("E18").Value = "Cofondateur & cogérant"
$ws.Range("E20").Value = "Médecin associé"
$ws.Range("E21").Value = "Cofondateur & cogérant"
$ws.Range("E22").Value = "Médecin associée"
$ws.Range("E23").Value = "Médecin associé"

# Update the view state: scroll the sheet and reselect E23 (best-effort;
# matches the author re-reviewing the last edited row).
$ws.Select()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("E23").Select()
